$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.569.49"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "3.363.36"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'562.41"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").Value = "'175.93"
$ws.Range("E6").Value = "  +3.81%  "

$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "  +4.15%  "

$ws.Range("D8").Value = "3.351.66"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.633"
$ws.Range("E10").Value = "  +4.75%  "

$ws.Range("D11").Value = "'0.164"
$ws.Range("E11").Value = "  +10.54%  "

$ws.Range("D12").Value = "'55.12"
$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "'0.0000275"
$ws.Range("E13").Value = "  +5.14%  "

$ws.Range("D14").Value = "'9.11"
$ws.Range("E14").Value = "  +3.38%  "

$ws.Range("D15").Value = "3.902.59"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "'18.29"
$ws.Range("E16").Value = "  +4.14%  "

$ws.Range("D17").Value = "3.363.28"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("E18").Value = "  +0.88%  "

$ws.Range("E19").Value = "  +3.10%  "

$ws.Range("D20").Value = "64.440.45"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").Value = "'0.992"
$ws.Range("E21").Value = "  +2.95%  "

$ws.Range("D22").Value = "'467.84"
$ws.Range("E22").Value = "  +18.10%  "

$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = "  +17.23%  "

$ws.Range("D24").Value = "'4.12"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("D25").Value = "'86.75"
$ws.Range("E25").Value = "  +6.82%  "

$ws.Range("D26").Value = "'13.51"
$ws.Range("E26").Value = "  +3.65%  "

$ws.Range("D27").Value = "'10.84"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").Value = "'2.84"
$ws.Range("E28").Value = "  +4.81%  "

$ws.Range("D29").Value = "'8.82"
$ws.Range("E29").Value = "  +3.40%  "

$ws.Range("D30").Value = "'30.22"
$ws.Range("E30").Value = "  +4.49%  "

$ws.Range("D31").Value = "'6.65"
$ws.Range("E31").Value = "  +3.55%  "

$ws.Range("D32").Value = "'11.50"

$ws.Range("D33").Value = "'579.40"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("E34").Value = "  +3.70%  "

$ws.Range("D35").Value = "'59.66"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -5.89%  "

$ws.Range("D38").Value = "'35.95"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "0.0₃0755"
$ws.Range("E39").Value = "  +4.03%  "

$ws.Range("D40").Value = "'3.45"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").Value = "'0.371"
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("D42").Value = "3.090.47"
$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").Value = "'2.82"
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").Value = "'2.52"
$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46").Value = "'0.0412"
$ws.Range("E46").Value = "  +3.58%  "

$ws.Range("D47").Value = "'3.20"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("E48").Value = "  +4.58%  "

$ws.Range("D49").Value = "'2.60"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'137.46"
$ws.Range("E50").Value = "  +4.29%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'8.38"
$ws.Range("E51").Value = "  +4.76%  "
